{"js": "// Replace the date line and each of the 25 division-fact table cells with\n// their new values. Every \"find\" string is unique in the document, so a\n// straightforward search + full-text replace on the matching range is safe.\nconst replacements = [\n  [\"2025-01-07 Tuesday\", \"2025-01-08 Wednesday\"],\n  [\"12\u00f76=2, 0\", \"97\u00f76=16, 1\"],\n  [\"18\u00f75=3, 3\", \"10\u00f72=5, 0\"],\n  [\"99\u00f73=33, 0\", \"38\u00f79=4, 2\"],\n  [\"13\u00f77=1, 6\", \"82\u00f78=10, 2\"],\n  [\"68\u00f78=8, 4\", \"68\u00f75=13, 3\"],\n  [\"45\u00f74=11, 1\", \"88\u00f75=17, 3\"],\n  [\"92\u00f76=15, 2\", \"86\u00f79=9, 5\"],\n  [\"75\u00f76=12, 3\", \"78\u00f75=15, 3\"],\n  [\"90\u00f72=45, 0\", \"88\u00f76=14, 4\"],\n  [\"96\u00f78=12, 0\", \"62\u00f75=12, 2\"],\n  [\"61\u00f73=20, 1\", \"14\u00f78=1, 6\"],\n  [\"35\u00f74=8, 3\", \"16\u00f74=4, 0\"],\n  [\"59\u00f78=7, 3\", \"55\u00f77=7, 6\"],\n  [\"43\u00f74=10, 3\", \"90\u00f76=15, 0\"],\n  [\"48\u00f73=16, 0\", \"47\u00f77=6, 5\"],\n  [\"30\u00f79=3, 3\", \"15\u00f72=7, 1\"],\n  [\"59\u00f76=9, 5\", \"64\u00f78=8, 0\"],\n  [\"39\u00f72=19, 1\", \"46\u00f76=7, 4\"],\n  [\"63\u00f76=10, 3\", \"80\u00f79=8, 8\"],\n  [\"66\u00f77=9, 3\", \"49\u00f79=5, 4\"],\n  [\"31\u00f73=10, 1\", \"82\u00f73=27, 1\"],\n  [\"66\u00f79=7, 3\", \"36\u00f77=5, 1\"],\n  [\"38\u00f76=6, 2\", \"65\u00f74=16, 1\"],\n  [\"31\u00f79=3, 4\", \"89\u00f77=12, 5\"],\n  [\"43\u00f76=7, 1\", \"86\u00f75=17, 1\"]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text: \" + oldText);\n  }\n\n  // Replace only the first (and, for this document, only) match, preserving\n  // the run's existing formatting by replacing the found range's text.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replace the date line and each of the 25 division-fact table cells with\n# their new values. Each \"find\" string occurs exactly once in the document,\n# so Find/Replace with wdReplaceOne is safe and keeps existing run formatting.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = '2025-01-07 Tuesday'; New = '2025-01-08 Wednesday' },\n    @{ Old = '12\u00f76=2, 0'; New = '97\u00f76=16, 1' },\n    @{ Old = '18\u00f75=3, 3'; New = '10\u00f72=5, 0' },\n    @{ Old = '99\u00f73=33, 0'; New = '38\u00f79=4, 2' },\n    @{ Old = '13\u00f77=1, 6'; New = '82\u00f78=10, 2' },\n    @{ Old = '68\u00f78=8, 4'; New = '68\u00f75=13, 3' },\n    @{ Old = '45\u00f74=11, 1'; New = '88\u00f75=17, 3' },\n    @{ Old = '92\u00f76=15, 2'; New = '86\u00f79=9, 5' },\n    @{ Old = '75\u00f76=12, 3'; New = '78\u00f75=15, 3' },\n    @{ Old = '90\u00f72=45, 0'; New = '88\u00f76=14, 4' },\n    @{ Old = '96\u00f78=12, 0'; New = '62\u00f75=12, 2' },\n    @{ Old = '61\u00f73=20, 1'; New = '14\u00f78=1, 6' },\n    @{ Old = '35\u00f74=8, 3'; New = '16\u00f74=4, 0' },\n    @{ Old = '59\u00f78=7, 3'; New = '55\u00f77=7, 6' },\n    @{ Old = '43\u00f74=10, 3'; New = '90\u00f76=15, 0' },\n    @{ Old = '48\u00f73=16, 0'; New = '47\u00f77=6, 5' },\n    @{ Old = '30\u00f79=3, 3'; New = '15\u00f72=7, 1' },\n    @{ Old = '59\u00f76=9, 5'; New = '64\u00f78=8, 0' },\n    @{ Old = '39\u00f72=19, 1'; New = '46\u00f76=7, 4' },\n    @{ Old = '63\u00f76=10, 3'; New = '80\u00f79=8, 8' },\n    @{ Old = '66\u00f77=9, 3'; New = '49\u00f79=5, 4' },\n    @{ Old = '31\u00f73=10, 1'; New = '82\u00f73=27, 1' },\n    @{ Old = '66\u00f79=7, 3'; New = '36\u00f77=5, 1' },\n    @{ Old = '38\u00f76=6, 2'; New = '65\u00f74=16, 1' },\n    @{ Old = '31\u00f79=3, 4'; New = '89\u00f77=12, 5' },\n    @{ Old = '43\u00f76=7, 1'; New = '86\u00f75=17, 1' }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $ok = $find.Execute($pair.Old, $true, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n    if (-not $ok) {\n        throw \"Could not find text: $($pair.Old)\"\n    }\n}\n"}
